$d = $word.ActiveDocument

# Remove the placeholder run "vnpt.SiteAddress" that follows "Địa chỉ: "
$d.Content.Find.Execute("vnpt.SiteAddress", $false, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)
